$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 56
$prev = $row - 1

# Copy formatting from the row above (matches existing per-column styling:
# bold/bordered/centered style on column A, date-number-format style on column E)
$ws.Cells.Item($prev, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item($prev, 5).Copy()
$ws.Cells.Item($row, 5).PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

$ws.Cells.Item($row, 1).Value = 55
$ws.Cells.Item($row, 2).Value = "azerbaijan"
$ws.Cells.Item($row, 3).Value = "premier-league"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45233.66666666666
$ws.Cells.Item($row, 6).Value = "Sabail"
$ws.Cells.Item($row, 7).Value = 2
$ws.Cells.Item($row, 8).Value = "Turan"
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 2.39
$ws.Cells.Item($row, 11).Value = "02/11/2023 06:42"
$ws.Cells.Item($row, 12).Value = 2.64
$ws.Cells.Item($row, 13).Value = "03/11/2023 12:42"
$ws.Cells.Item($row, 14).Value = 3.05
$ws.Cells.Item($row, 15).Value = "02/11/2023 06:42"
$ws.Cells.Item($row, 16).Value = 3.25
$ws.Cells.Item($row, 17).Value = "03/11/2023 14:05"
$ws.Cells.Item($row, 18).Value = 2.79
$ws.Cells.Item($row, 19).Value = "02/11/2023 06:42"
$ws.Cells.Item($row, 20).Value = 2.56
$ws.Cells.Item($row, 21).Value = "03/11/2023 12:47"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/azerbaijan/premier-league/sabail-turan/42d5V9Li/"
